$d = $word.ActiveDocument

# Find the paragraph holding the "Ver no Jupiter ..." site-footer line and the one
# holding the "(c) 2020 ... Contact: luizeleno@usp.br ..." copyright line.
$jupiterPara = $null
$copyrightPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $jupiterPara = $p
    }
    elseif ($t -like "*Contact: luizeleno@usp.br*") {
        $copyrightPara = $p
    }
}

# The blank separator paragraph right before "Ver no Jupiter ..." is removed together
# with it and the copyright paragraph that follows, collapsing the pair of blank
# paragraphs that used to sandwich this footer block down to the single one that
# remains right after "Geophysics. Springer - Verlag, Berlin, 171p."
$blankBefore = $jupiterPara.Previous(1)

$start = $blankBefore.Range.Start
$end = $copyrightPara.Range.End

$r = $d.Range($start, $end)
$r.Delete()
